# history.xlsx — "done gui for web"
#
# The sheet had two date-grouped rows (A2:A3 merged under "30-03-2018",
# A4:A5 merged under "23-02-2018"), each with a "12"/"23" pair in columns
# B/C for the first group and a "21"/"41" pair for the second, all text.
#
# The edit keeps only the first group, updates its date and its second
# figure, and drops column C entirely along with the second group/date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date in A2 and the count in B3, forcing them to stay text
# (these values look numeric/date-like, so Excel would otherwise coerce
# them on assignment).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "08-05-2018"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2"

# Drop column C (the "23"/"41" values) for the remaining rows.
$ws.Range("C2:C3").Clear()

# Remove the second date group (rows 4-5) entirely - this also drops the
# A4:A5 merge and shrinks the used range down to A2:B3.
$ws.Rows("4:5").Delete()
